$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped from
# 45190 (2023-09-21) to 45192 (2023-09-23) for every data row (rows 2-556).
$ws.Range("C2:C556").Value = 45192
